$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 441.53845
$ws.Range("J2").Value = 533
$ws.Range("L2").Value = 533
$ws.Range("N2").Value = -759

# Row 9
$ws.Range("H9").Value = 1142299.6
$ws.Range("I9").Value = 516
$ws.Range("J9").Value = 1499107
$ws.Range("K9").Value = 516
$ws.Range("L9").Value = 1499107
$ws.Range("M9").Value = -347
$ws.Range("N9").Value = -1499445

# Row 11
$ws.Range("H11").Value = 4860.091
$ws.Range("I11").Value = 4860.091
$ws.Range("K11").Value = 4860.091
$ws.Range("M11").Value = -4720.091

# Row 17
$ws.Range("H17").Value = 3642.842
$ws.Range("J17").Value = 3789.6667
$ws.Range("L17").Value = 11369.0001
$ws.Range("N17").Value = -11705.0001

# Row 18
$ws.Range("H18").Value = 2493.3333
$ws.Range("I18").Value = 2493.3333
$ws.Range("K18").Value = 2493.3333
$ws.Range("M18").Value = -2209.3333

# Row 40
$ws.Range("H40").Value = 3049.375
$ws.Range("I40").Value = 3825.5
$ws.Range("J40").Value = 2273.25
$ws.Range("K40").Value = 3825.5
$ws.Range("L40").Value = 2273.25
$ws.Range("M40").Value = -3650.5
$ws.Range("N40").Value = -2623.25

# Row 43
$ws.Range("H43").Value = 3100.1667
$ws.Range("I43").Value = 2633.3333
$ws.Range("K43").Value = 2633.3333
$ws.Range("M43").Value = -2564.3333

# Row 69
$ws.Range("H69").Value = 12660.333
$ws.Range("I69").Value = 20982.5
$ws.Range("K69").Value = 62947.5
$ws.Range("M69").Value = -62073.5

# Row 72
$ws.Range("H72").Value = 12660.333
$ws.Range("I72").Value = 20982.5
$ws.Range("K72").Value = 188842.5
$ws.Range("M72").Value = -184474.5

# Row 74
$ws.Range("H74").Value = 7120.625
$ws.Range("I74").Value = 6157
$ws.Range("J74").Value = 7698.8
$ws.Range("K74").Value = 6157
$ws.Range("L74").Value = 7698.8
$ws.Range("M74").Value = -5221
$ws.Range("N74").Value = -9570.799999999999

# Row 77
$ws.Range("H77").Value = 7120.625
$ws.Range("I77").Value = 6157
$ws.Range("J77").Value = 7698.8
$ws.Range("K77").Value = 30785
$ws.Range("L77").Value = 38494
$ws.Range("M77").Value = -26105
$ws.Range("N77").Value = -47854

# Row 137
$ws.Range("H137").Value = 3230.425
$ws.Range("I137").Value = 2412.0454
$ws.Range("J137").Value = 4230.6665
$ws.Range("K137").Value = 7236.1362
$ws.Range("L137").Value = 12691.9995
$ws.Range("M137").Value = -4686.1362
$ws.Range("N137").Value = -17791.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2387.3333
$ws.Range("I45").Value = 2064.8
$ws.Range("K45").Value = 2064.8
$ws.Range("M45").Value = -1687.8

# Row 74
$ws.Range("H74").Value = 6216.387
$ws.Range("I74").Value = 5025.75
$ws.Range("K74").Value = 5025.75
$ws.Range("M74").Value = -4151.75

# Row 77
$ws.Range("H77").Value = 6216.387
$ws.Range("I77").Value = 5025.75
$ws.Range("K77").Value = 25128.75
$ws.Range("M77").Value = -20760.75

# Row 102
$ws.Range("H102").Value = 71572460
$ws.Range("I102").Value = 1814.75
$ws.Range("K102").Value = 1814.75
$ws.Range("M102").Value = -192.75

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 52437.547
$ws.Range("I82").Value = 34614.5
$ws.Range("K82").Value = 34614.5
$ws.Range("M82").Value = -34231.5

# Row 85
$ws.Range("H85").Value = 52437.547
$ws.Range("I85").Value = 34614.5
$ws.Range("K85").Value = 34614.5
$ws.Range("M85").Value = -33288.5

# Row 86
$ws.Range("H86").Value = 335152
$ws.Range("I86").Value = 1706
$ws.Range("K86").Value = 1706
$ws.Range("M86").Value = -583

# Row 88
$ws.Range("H88").Value = 34147.438
$ws.Range("J88").Value = 34147.438
$ws.Range("L88").Value = 34147.438
$ws.Range("N88").Value = -34959.438

# Row 89
$ws.Range("H89").Value = 335152
$ws.Range("I89").Value = 1706
$ws.Range("K89").Value = 8530
$ws.Range("M89").Value = -2914

# Row 91
$ws.Range("H91").Value = 34147.438
$ws.Range("J91").Value = 34147.438
$ws.Range("L91").Value = 34147.438
$ws.Range("N91").Value = -36955.438

# Row 99
$ws.Range("H99").Value = 3566.6667
$ws.Range("I99").Value = 3350
$ws.Range("K99").Value = 3350
$ws.Range("M99").Value = -1852

# Row 125
$ws.Range("H125").Value = 39139.75
$ws.Range("J125").Value = 39139.75
$ws.Range("L125").Value = 39139.75
$ws.Range("N125").Value = -48979.75

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 31
$ws.Range("H31").Value = 3939.9
$ws.Range("I31").Value = 2249.7222
$ws.Range("K31").Value = 2249.7222
$ws.Range("M31").Value = -1954.7222

# Row 34
$ws.Range("H34").Value = 3939.9
$ws.Range("I34").Value = 2249.7222
$ws.Range("K34").Value = 2249.7222
$ws.Range("M34").Value = -2047.7222

# Row 62
$ws.Range("H62").Value = 10638.077
$ws.Range("I62").Value = 8366.111000000001
$ws.Range("K62").Value = 8366.111000000001
$ws.Range("M62").Value = -7742.111000000001

# Row 65
$ws.Range("H65").Value = 10638.077
$ws.Range("I65").Value = 8366.111000000001
$ws.Range("K65").Value = 41830.55500000001
$ws.Range("M65").Value = -38710.55500000001

# Row 68
$ws.Range("H68").Value = 391798.8
$ws.Range("J68").Value = 386331.34
$ws.Range("L68").Value = 386331.34
$ws.Range("N68").Value = -387829.34

# Row 71
$ws.Range("H71").Value = 391798.8
$ws.Range("J71").Value = 386331.34
$ws.Range("L71").Value = 1158994.02
$ws.Range("N71").Value = -1166482.02

# Row 109
$ws.Range("H109").Value = 66704600
$ws.Range("J109").Value = 66704600
$ws.Range("L109").Value = 66704600
$ws.Range("N109").Value = -66706680

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 53232.367
$ws.Range("I44").Value = 701.1429000000001
$ws.Range("J44").Value = 200319.8
$ws.Range("K44").Value = 2103.4287
$ws.Range("L44").Value = 600959.3999999999
$ws.Range("M44").Value = -1705.4287
$ws.Range("N44").Value = -601755.3999999999

# Row 132
$ws.Range("H132").Value = 81665.69500000001
$ws.Range("I132").Value = 206941.2
$ws.Range("K132").Value = 1862470.8
$ws.Range("M132").Value = -1859940.8

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888
$ws.Range("N2").ClearContents()

# Row 46
$ws.Range("H46").Value = 7243.2
$ws.Range("I46").Value = 1933.3334
$ws.Range("K46").Value = 1933.3334
$ws.Range("M46").Value = -1745.3334

# Row 64
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5450

# Row 67
$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6560

# Row 68
$ws.Range("H68").Value = 4238.778
$ws.Range("J68").Value = 6150
$ws.Range("L68").Value = 6150
$ws.Range("N68").Value = -7648

# Row 71
$ws.Range("H71").Value = 4238.778
$ws.Range("J71").Value = 6150
$ws.Range("L71").Value = 30750
$ws.Range("N71").Value = -38238

# Row 82
$ws.Range("H82").Value = 1985.2858
$ws.Range("I82").Value = 1874.25
$ws.Range("K82").Value = 1874.25
$ws.Range("M82").Value = -1513.25

# Row 85
$ws.Range("H85").Value = 1985.2858
$ws.Range("I85").Value = 1874.25
$ws.Range("K85").Value = 1874.25
$ws.Range("M85").Value = -626.25

# Row 136
$ws.Range("H136").Value = 1908233.6
$ws.Range("I136").Value = 2471488.2
$ws.Range("K136").Value = 7414464.600000001
$ws.Range("M136").Value = -7411914.600000001

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 233229.56
$ws.Range("J2").Value = 339333.34
$ws.Range("L2").Value = 339333.34
$ws.Range("N2").Value = -339557.34

# Row 62
$ws.Range("H62").Value = 9018.182000000001
$ws.Range("I62").Value = 8224.25
$ws.Range("J62").Value = 9471.857
$ws.Range("K62").Value = 8224.25
$ws.Range("L62").Value = 9471.857
$ws.Range("M62").Value = -7600.25
$ws.Range("N62").Value = -10719.857

# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 65
$ws.Range("H65").Value = 9018.182000000001
$ws.Range("I65").Value = 8224.25
$ws.Range("J65").Value = 9471.857
$ws.Range("K65").Value = 41121.25
$ws.Range("L65").Value = 47359.285
$ws.Range("M65").Value = -38001.25
$ws.Range("N65").Value = -53599.285

# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 81
$ws.Range("H81").Value = 13074.389
$ws.Range("I81").Value = 6127.143
$ws.Range("K81").Value = 12254.286
$ws.Range("M81").Value = -11193.286

# Row 84
$ws.Range("H84").Value = 13074.389
$ws.Range("I84").Value = 6127.143
$ws.Range("K84").Value = 61271.43
$ws.Range("M84").Value = -55967.43

# Row 109
$ws.Range("H109").Value = 56394.5
$ws.Range("J109").Value = 56394.5
$ws.Range("L109").Value = 56394.5
$ws.Range("N109").Value = -59168.5
